$d = $word.ActiveDocument

$replacements = @(
    @{ old = "2025-02-19 Wednesday"; new = "2025-02-20 Thursday" },
    @{ old = "26×14=364";  new = "65×47=3055" },
    @{ old = "53×43=2279"; new = "78×64=4992" },
    @{ old = "25×34=850";  new = "99×63=6237" },
    @{ old = "21×42=882";  new = "96×51=4896" },
    @{ old = "74×97=7178"; new = "37×56=2072" },
    @{ old = "83×87=7221"; new = "80×40=3200" },
    @{ old = "85×12=1020"; new = "97×97=9409" },
    @{ old = "27×75=2025"; new = "48×17=816"  },
    @{ old = "15×87=1305"; new = "39×72=2808" },
    @{ old = "22×97=2134"; new = "15×85=1275" },
    @{ old = "34×40=1360"; new = "36×11=396"  },
    @{ old = "27×19=513";  new = "85×29=2465" },
    @{ old = "14×32=448";  new = "42×50=2100" },
    @{ old = "92×42=3864"; new = "97×67=6499" },
    @{ old = "29×50=1450"; new = "32×40=1280" },
    @{ old = "41×66=2706"; new = "88×23=2024" },
    @{ old = "45×21=945";  new = "97×42=4074" },
    @{ old = "65×91=5915"; new = "15×32=480"  },
    @{ old = "97×32=3104"; new = "43×12=516"  },
    @{ old = "76×35=2660"; new = "50×70=3500" },
    @{ old = "84×21=1764"; new = "50×62=3100" },
    @{ old = "53×27=1431"; new = "78×45=3510" },
    @{ old = "75×23=1725"; new = "81×77=6237" },
    @{ old = "94×85=7990"; new = "15×50=750"  },
    @{ old = "60×13=780";  new = "42×46=1932" }
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
